# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (so the sheet
#    order becomes 总计, 2022-Q4, 2022-Q2, 2021-Q2) and populate it with
#    the fund-holding breakdown for the quarter.
# 2) Update the "总计" (summary) sheet: the existing "2022-Q2" row now
#    carries the new "2022-Q4" numbers, and the old rows shift down to
#    make room, ending with 2022-Q4 / 2022-Q2 / 2021-Q2 in that order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计" (index 1).
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $zongji)
$q4.Name = "2022-Q4"

# Header row.
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Columns B-G hold numeric-looking text (fund codes / percentages) that
# must stay text, not auto-convert to numbers - force Text format first.
$q4.Range("B2:G10").NumberFormat = "@"

# Row 2 - 513330
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "513330"
$q4.Cells.Item(2,3).Value = "华夏恒生互联网科技业ETF（QDII）"
$q4.Cells.Item(2,4).Value = "246.27"
$q4.Cells.Item(2,5).Value = "96.18"
$q4.Cells.Item(2,6).Value = "3.29"
$q4.Cells.Item(2,7).Value = "8.1023"
$q4.Cells.Item(2,8).Value = 9

# Row 3 - 012208
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "012208"
$q4.Cells.Item(3,3).Value = "华夏港股前沿经济混合（QDII）A"
$q4.Cells.Item(3,4).Value = "10.38"
$q4.Cells.Item(3,5).Value = "92.91"
$q4.Cells.Item(3,6).Value = "4.14"
$q4.Cells.Item(3,7).Value = "0.4297"
$q4.Cells.Item(3,8).Value = 8

# Row 4 - 003243
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "003243"
$q4.Cells.Item(4,3).Value = "上投摩根中国世纪灵活配置混合人民币份额（QDII）"
$q4.Cells.Item(4,4).Value = "1.24"
$q4.Cells.Item(4,5).Value = "85.53"
$q4.Cells.Item(4,6).Value = "4.47"
$q4.Cells.Item(4,7).Value = "0.0554"
$q4.Cells.Item(4,8).Value = 3

# Row 5 - 003244
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "003244"
$q4.Cells.Item(5,3).Value = "上投摩根中国世纪灵活配置混合美元现钞（QDII）"
$q4.Cells.Item(5,4).Value = "1.24"
$q4.Cells.Item(5,5).Value = "85.53"
$q4.Cells.Item(5,6).Value = "4.47"
$q4.Cells.Item(5,7).Value = "0.0554"
$q4.Cells.Item(5,8).Value = 3

# Row 6 - 003245
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "003245"
$q4.Cells.Item(6,3).Value = "上投摩根中国世纪灵活配置混合美元现汇（QDII）"
$q4.Cells.Item(6,4).Value = "1.24"
$q4.Cells.Item(6,5).Value = "85.53"
$q4.Cells.Item(6,6).Value = "4.47"
$q4.Cells.Item(6,7).Value = "0.0554"
$q4.Cells.Item(6,8).Value = 3

# Row 7 - 460010
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "460010"
$q4.Cells.Item(7,3).Value = "华泰柏瑞亚洲领导企业混合（QDII）"
$q4.Cells.Item(7,4).Value = "0.52"
$q4.Cells.Item(7,5).Value = "97.17"
$q4.Cells.Item(7,6).Value = "4.85"
$q4.Cells.Item(7,7).Value = "0.0252"
$q4.Cells.Item(7,8).Value = 10

# Row 8 - 015884
$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).Value = "015884"
$q4.Cells.Item(8,3).Value = "中欧港股数字经济混合（QDII）A"
$q4.Cells.Item(8,4).Value = "0.21"
$q4.Cells.Item(8,5).Value = "88.32"
$q4.Cells.Item(8,6).Value = "7.53"
$q4.Cells.Item(8,7).Value = "0.0158"
$q4.Cells.Item(8,8).Value = 5

# Row 9 - 012209
$q4.Cells.Item(9,1).Value = 7
$q4.Cells.Item(9,2).Value = "012209"
$q4.Cells.Item(9,3).Value = "华夏港股前沿经济混合（QDII）C"
$q4.Cells.Item(9,4).Value = "0.38"
$q4.Cells.Item(9,5).Value = "92.91"
$q4.Cells.Item(9,6).Value = "4.14"
$q4.Cells.Item(9,7).Value = "0.0157"
$q4.Cells.Item(9,8).Value = 8

# Row 10 - 015885
$q4.Cells.Item(10,1).Value = 8
$q4.Cells.Item(10,2).Value = "015885"
$q4.Cells.Item(10,3).Value = "中欧港股数字经济混合（QDII）C"
$q4.Cells.Item(10,4).Value = "0.12"
$q4.Cells.Item(10,5).Value = "88.32"
$q4.Cells.Item(10,6).Value = "7.53"
$q4.Cells.Item(10,7).Value = "0.0090"
$q4.Cells.Item(10,8).Value = 5

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: row 2 becomes 2022-Q4, the old
#    2022-Q2 row moves to row 3, and 2021-Q2 moves to row 4.
# ---------------------------------------------------------------------

# Row 4 first (uses the old row-3 style for the new A4 cell) before the
# old row-3 values are overwritten.
$zongji.Cells.Item(3,1).Copy($zongji.Cells.Item(4,1))
$zongji.Cells.Item(4,1).Value = 2
$zongji.Cells.Item(4,2).Value = "2021-Q2"
$zongji.Cells.Item(4,3).Value = 6
$zongji.Cells.Item(4,4).Value = 4.82

# Row 3 becomes the former 2022-Q2 entry.
$zongji.Cells.Item(3,1).Value = 1
$zongji.Cells.Item(3,2).Value = "2022-Q2"
$zongji.Cells.Item(3,3).Value = 5
$zongji.Cells.Item(3,4).Value = 1.83

# Row 2 becomes the new 2022-Q4 entry.
$zongji.Cells.Item(2,2).Value = "2022-Q4"
$zongji.Cells.Item(2,3).Value = 9
$zongji.Cells.Item(2,4).Value = 8.76

# ---------------------------------------------------------------------
# 3. Restore the original active sheet (the last sheet, 2021-Q2, was
#    the one selected before this edit).
# ---------------------------------------------------------------------
$lastws = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastws.Activate()
